$d = $word.ActiveDocument

# Locate the "Experienced DevOps Engineer ..." profile-summary paragraph,
# and the (empty) paragraph that immediately precedes it, by content
# rather than a hard-coded index, so the script is resilient to minor
# paragraph-numbering differences.
$profileIdx = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "Experienced DevOps Engineer and Frontend Developer*") {
        $profileIdx = $i
    }
}
if ($profileIdx -eq -1) {
    throw "Could not locate the profile-summary paragraph"
}
$bookmarkIdx = $profileIdx - 1

# --- Change 1: add the _GoBack bookmark to the empty paragraph that
#     precedes the "Experienced..." paragraph. ---
$pBookmark = $d.Paragraphs($bookmarkIdx)
$rBookmark = $pBookmark.Range
$rBookmark.Collapse(1) | Out-Null
$rBookmark.InsertBefore("X")
$d.Bookmarks.Add("_GoBack", $rBookmark) | Out-Null
$rBookmarkChar = $d.Range($rBookmark.Start, $rBookmark.Start + 1)
$rBookmarkChar.Delete() | Out-Null

# --- Change 2: replace the whole "Experienced ... deployment speed."
#     paragraph body with a single unformatted run containing the new
#     profile summary text. ---
$pProfile = $d.Paragraphs($profileIdx)
$rProfile = $pProfile.Range
$rProfile.MoveEnd(1, -1) | Out-Null
$rProfile.Text = ""

$pProfile2 = $d.Paragraphs($profileIdx)
$rProfile2 = $pProfile2.Range
$rProfile2.MoveEnd(1, -1) | Out-Null
$newText = "Experienced DevOps Engineer and Frontend Developer with over 5 years" + [char]8217 + " experience and a strong background in cloud platforms like Jenkins, Kubernetes, GitOps, EKS, RDS, Terraform, AWS, SQL, MongoDB, Python, Groovy, Spring Boot, .NET Core, Argo Rollouts, Helm, CI/CD, Velero, Istio, Lua, CI/CD pipelines including Jenkins, Kubernetes, GitOps, EKS, RDS, Terraform, AWS, SQL, MongoDB, Python, Groovy, Spring Boot, .NET Core, Argo Rollouts, Helm, CI/CD, Velero, Istio, Lua, and frameworks such as Jenkins, Kubernetes, GitOps, EKS, RDS, Terraform, AWS, SQL, MongoDB, Python, Groovy, Spring Boot, .NET Core, Argo Rollouts, Helm, CI/CD, Velero, Istio, Lua., "
$rProfile2.InsertAfter($newText)

Write-Host "Edit complete"
